$d = $word.ActiveDocument

# Locate the existing sentence and append the new sentence's text right
# after it (still inside the same run/paragraph as the original text, so
# the trailing bookmark stays glued to the end of the combined text for
# now).
$rng = $d.Content
$rng.Find.Execute("This word document is for git practice") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("This practice will be used in real-world projects")

# Now split the paragraph in two right before the newly-added sentence,
# turning it into its own paragraph. Splitting here (rather than exactly
# at the old end-of-text position) carries the trailing bookmark along
# with the text that follows the split point, so it ends up anchored to
# the new second paragraph - matching a normal "type new text after the
# old text, then press Enter" edit.
$rng2 = $d.Content
$rng2.Find.Execute("This practice will be used in real-world projects") | Out-Null
$rng2.Collapse(1)
$rng2.InsertParagraphBefore() | Out-Null
